$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (Planvale Services): tradeLocation H4 and tradeNameLocation D4 updated
$ws.Range("H4").Value = "Holyhead"
$ws.Range("D4").Value = "Shrewsbury"

# Row 6 (Autopics Garage Services): tradeLocation H6 updated
$ws.Range("H6").Value = "Rustington"

# Column C was widened (tradeName column)
$ws.Columns.Item(3).ColumnWidth = 25.6

# Selection moved to G6
$ws.Range("G6").Select()
